$d = $word.ActiveDocument

# Locate the paragraph that ends with "Αναφερόμενους εκπαιδευτικούς"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a", "`n") -eq "Αναφερόμενους εκπαιδευτικούς") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Target paragraph not found"
}

# Insert a new paragraph right after the target paragraph
$newRange = $target.Range.InsertParagraphAfter()

# The newly created paragraph is the one following $target
$newPara = $target.Next()
$newPara.Range.Text = "(μέσω της σχολικής μονάδας)"

$newPara.Style = $d.Styles.Item("a8")
$newPara.Range.ParagraphFormat.LeftIndent = 540
$newPara.Range.ParagraphFormat.FirstLineIndent = 360

$r = $newPara.Range
$r.Font.Name = "Calibri"
$r.Font.NameFarEast = "MS Mincho;ＭＳ 明朝"
$r.Font.NameAscii = "Calibri"
$r.Font.NameOther = "Times New Roman"
$r.Font.Size = 11
